$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "58.716.94"
$ws.Range("E2").Value = "  -2.61%  "
Set-TextValue "D3" "2.656.85"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "523.85"
$ws.Range("E5").Value = "  +0.16%  "
Set-TextValue "D6" "143.94"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -1.09%  "
Set-TextValue "D9" "6.91"
$ws.Range("E9").Value = "  +6.84%  "
$ws.Range("E10").Value = "  -3.40%  "
Set-TextValue "D11" "0.335"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  +1.37%  "
Set-TextValue "D13" "3.126.29"
$ws.Range("E13").Value = "  -0.86%  "
Set-TextValue "D14" "58.742.85"
$ws.Range("E14").Value = "  -2.72%  "
Set-TextValue "D15" "20.98"
$ws.Range("E16").Value = "  -1.67%  "
Set-TextValue "D17" "2.665.26"
$ws.Range("E17").Value = "  -1.20%  "
Set-TextValue "D18" "338.83"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("E20").Value = "  -1.51%  "
Set-TextValue "D21" "6.43"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E22").Value = "  -0.11%  "
Set-TextValue "D23" "63.93"
$ws.Range("E23").Value = "  +1.10%  "
Set-TextValue "D24" "0.424"
$ws.Range("E24").Value = "  +0.36%  "
Set-TextValue "D25" "0.166"
$ws.Range("E25").Value = "  -1.59%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  -1.63%  "
Set-TextValue "D28" "7.14"
$ws.Range("E28").Value = "  -2.91%  "
Set-TextValue "D29" "6.69"
$ws.Range("E29").Value = "  -1.30%  "
Set-TextValue "D30" "0.998"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -0.01%  "
Set-TextValue "D32" "18.89"
$ws.Range("E32").Value = "  -1.41%  "
Set-TextValue "D33" "151.01"
$ws.Range("E33").Value = "  +2.54%  "
Set-TextValue "D34" "4.17"
$ws.Range("E34").Value = "  -3.24%  "
Set-TextValue "D35" "0.933"
$ws.Range("E35").Value = "  -2.26%  "
Set-TextValue "D36" "1.18"
$ws.Range("E36").Value = "  -5.83%  "
Set-TextValue "D37" "0.877"
$ws.Range("E37").Value = "  -0.28%  "
Set-TextValue "D38" "36.78"
$ws.Range("E38").Value = "  -0.44%  "
Set-TextValue "D39" "1.44"
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D41" "0.999"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.610"
$ws.Range("E42").Value = "  -0.01%  "
Set-TextValue "D43" "275.99"
$ws.Range("E43").Value = "  -2.87%  "
Set-TextValue "D44" "19.69"
$ws.Range("E44").Value = "  -1.85%  "
Set-TextValue "D45" "0.0967"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0535"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D48" "2.061.20"
$ws.Range("E48").Value = "  -3.40%  "
Set-TextValue "D49" "4.73"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("E50").Value = "  -3.06%  "
Set-TextValue "D51" "18.69"
$ws.Range("E51").Value = "  -3.84%  "
